# Add two more soccer teams to the SoccerPage sheet.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A11").Value = "ARSENAL"
$ws1.Range("A12").Value = "TOTTENHAM HOTSPUR"
$null = $ws1.Range("A13").Select()

# Add Paris Saint-Germain stat-page rows to the TeamInfo sheet.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns.Item(1).ColumnWidth = 33.36328125
$ws2.Range("A8").Value = "Paris Saint-Germain Performance Stats"
$ws2.Range("A9").Value = "Paris Saint-Germain Discipline Stats"
$ws2.Range("A10").Value = "Paris Saint-Germain Scoring Stats"
$ws2.Range("A11").Value = "Paris Saint-Germain Results"
$ws2.Range("A11").WrapText = $true
$ws2.Range("A12").Value = "Paris Saint-Germain Squad"
$ws2.Range("A13").Value = "Paris Saint-Germain Transfers"

# Add a brand new MMA sheet with a UFC champions assertion row.
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "MMA"
$ws3.Range("A1").Value = "Assertions"
$ws3.Range("A2").Value = "Current and all-time UFC champions"
$ws3.Columns.Item(1).ColumnWidth = 31.26953125
$null = $ws3.Activate()
